$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.875.64'
$ws.Range("E2").Value = '  +0.55%  '

$ws.Range("D3").Value = '1.660.89'
$ws.Range("E3").Value = '  -1.89%  '

$ws.Range("D4").Value = '''1.001'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.44%  '

$ws.Range("D5").Value = '''316.57'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +2.11%  '

$ws.Range("D6").Value = '''0.9956'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.52%  '

$ws.Range("D7").Value = '''0.3634'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -2.11%  '

$ws.Range("D8").Value = '''47.19'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -3.50%  '

$ws.Range("D9").Value = '''0.3270'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -3.85%  '

$ws.Range("D10").Value = '''1.141'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -4.82%  '

$ws.Range("D11").Value = '''0.07070'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -4.66%  '

$ws.Range("D12").Value = '''0.9984'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.26%  '

$ws.Range("D13").Value = '''6.043'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -3.61%  '

$ws.Range("D14").Value = '''19.62'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -5.60%  '

$ws.Range("D15").Value = '1.669.02'
$ws.Range("E15").Value = '  -1.24%  '

$ws.Range("D16").Value = '''6.628'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -4.20%  '

$ws.Range("D17").Value = '''0.00001050'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -5.83%  '

$ws.Range("D18").Value = '''0.06615'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -1.17%  '

$ws.Range("D19").Value = '''0.9956'
$ws.Range("D19").ClearFormats()

$ws.Range("D20").Value = '''79.32'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -4.18%  '

$ws.Range("D21").Value = '''5.926'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -5.55%  '

$ws.Range("D22").Value = '''15.76'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -7.52%  '

$ws.Range("D23").Value = '''12.55'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -1.85%  '

$ws.Range("D24").Value = '24.911.08'
$ws.Range("E24").Value = '  +0.75%  '

$ws.Range("D25").Value = '''2.431'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.51%  '

$ws.Range("D26").Value = '''2.398'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -12.38%  '

$ws.Range("D27").Value = '''148.52'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.22%  '

$ws.Range("D28").Value = '''18.65'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -7.13%  '

$ws.Range("B29").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C29").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D29").Value = '1.851.77'
$ws.Range("E29").Value = '  -1.32%  '

$ws.Range("B30").Value = 'ImmutableX'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D30").Value = '''1.223'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -0.73%  '

$ws.Range("D31").Value = '''125.77'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -4.20%  '

$ws.Range("D32").Value = '''4.112'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -2.30%  '

$ws.Range("D33").Value = '''5.841'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -12.24%  '

$ws.Range("D34").Value = '''0.08448'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -2.37%  '

$ws.Range("D35").Value = '''1.677'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -4.07%  '

$ws.Range("D36").Value = '''12.27'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -9.17%  '

$ws.Range("D37").Value = '''1.286'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +2.76%  '

$ws.Range("D38").Value = '''5.201'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -5.03%  '

$ws.Range("D39").Value = '''0.02238'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -6.05%  '

$ws.Range("D40").Value = '''0.06037'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -7.98%  '

$ws.Range("D41").Value = '''0.2067'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -5.46%  '

$ws.Range("D42").Value = '''8.260'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -8.32%  '

$ws.Range("D43").Value = '''0.9943'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -0.63%  '

$ws.Range("D44").Value = '''0.5937'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -6.40%  '

$ws.Range("D45").Value = '''3.817'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.38%  '

$ws.Range("D46").Value = '''12.84'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -5.51%  '

$ws.Range("D47").Value = '''0.5649'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -6.05%  '

$ws.Range("D48").Value = '''125.31'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -1.74%  '

$ws.Range("D49").Value = '''1.954'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -6.40%  '

$ws.Range("D50").Value = '''0.07028'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -2.63%  '

$ws.Range("D51").Value = '''1.195'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.95%  '
